$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 31   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/1/2024  Through  7/7/2024"

# --- Crime statistics table updates (rows 14-33) ---
$ws.Range("C14").Value = 2
$ws.Range("F14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("G14").Value = "'0"
$ws.Range("E14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("I14").Value = 19
$ws.Range("K14").Value = 111.111111111111
$ws.Range("L14").Value = -24
$ws.Range("M14").Value = -42.424242424242
$ws.Range("N14").Value = -73.239436619718

$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 500
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 112.5
$ws.Range("I15").Value = 93
$ws.Range("J15").Value = 74
$ws.Range("K15").Value = 25.675675675675
$ws.Range("L15").Value = -10.576923076923
$ws.Range("M15").Value = 30.985915492957
$ws.Range("N15").Value = -52.061855670103

$ws.Range("C16").Value = 22
$ws.Range("D16").Value = 29
$ws.Range("E16").Value = -24.137931034482
$ws.Range("F16").Value = 108
$ws.Range("G16").Value = 103
$ws.Range("H16").Value = 4.854368932038
$ws.Range("I16").Value = 696
$ws.Range("J16").Value = 657
$ws.Range("K16").Value = 5.93607305936
$ws.Range("L16").Value = -2.793296089385
$ws.Range("M16").Value = -34.953271028037
$ws.Range("N16").Value = -81.936153646509

$ws.Range("C17").Value = 62
$ws.Range("D17").Value = 77
$ws.Range("E17").Value = -19.480519480519
$ws.Range("G17").Value = 275
$ws.Range("H17").Value = -3.272727272727
$ws.Range("I17").Value = 1609
$ws.Range("J17").Value = 1437
$ws.Range("K17").Value = 11.96938065414
$ws.Range("L17").Value = 15.58908045977
$ws.Range("M17").Value = 88.849765258216
$ws.Range("N17").Value = -11.980306345733

$ws.Range("C18").Value = 18
$ws.Range("D18").Value = 21
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 68
$ws.Range("G18").Value = 83
$ws.Range("H18").Value = -18.072289156626
$ws.Range("I18").Value = 511
$ws.Range("J18").Value = 574
$ws.Range("K18").Value = -10.975609756097
$ws.Range("L18").Value = -8.258527827648
$ws.Range("M18").Value = -51.006711409396
$ws.Range("N18").Value = -89.067180145485

$ws.Range("C19").Value = 66
$ws.Range("D19").Value = 65
$ws.Range("E19").Value = 1.538461538461
$ws.Range("F19").Value = 226
$ws.Range("G19").Value = 271
$ws.Range("H19").Value = -16.60516605166
$ws.Range("I19").Value = 1643
$ws.Range("J19").Value = 1806
$ws.Range("K19").Value = -9.025470653377
$ws.Range("L19").Value = -14.516129032258
$ws.Range("M19").Value = 16.524822695035
$ws.Range("N19").Value = -59.946367625548

$ws.Range("C20").Value = 34
$ws.Range("D20").Value = 29
$ws.Range("E20").Value = 17.241379310344
$ws.Range("F20").Value = 149
$ws.Range("G20").Value = 159
$ws.Range("H20").Value = -6.2893081761
$ws.Range("I20").Value = 949
$ws.Range("J20").Value = 885
$ws.Range("K20").Value = 7.231638418079
$ws.Range("L20").Value = 18.625
$ws.Range("M20").Value = 13.516746411483
$ws.Range("N20").Value = -90.281618023553

$ws.Range("C21").Value = 210
$ws.Range("D21").Value = 222
$ws.Range("E21").Value = -5.405405405405
$ws.Range("F21").Value = 837
$ws.Range("G21").Value = 899
$ws.Range("H21").Value = -6.896551724137
$ws.Range("I21").Value = 5520
$ws.Range("J21").Value = 5442
$ws.Range("K21").Value = 1.433296582138
$ws.Range("L21").Value = 0.072516316171
$ws.Range("M21").Value = 3.857008466603
$ws.Range("N21").Value = -77.457426389512

$ws.Range("C22").Value = 4
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -20
$ws.Range("I22").Value = 65
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 8.333333333333
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = 4.838709677419

$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = 5.263157894736
$ws.Range("I23").Value = 122
$ws.Range("J23").Value = 122
$ws.Range("L23").Value = 6.086956521739
$ws.Range("M23").Value = 50.617283950617

$ws.Range("C24").Value = 174
$ws.Range("E24").Value = -2.247191011235
$ws.Range("F24").Value = 740
$ws.Range("G24").Value = 755
$ws.Range("H24").Value = -1.986754966887
$ws.Range("I24").Value = 4642
$ws.Range("J24").Value = 4769
$ws.Range("K24").Value = -2.663032082197
$ws.Range("L24").Value = -7.16
$ws.Range("M24").Value = 49.164524421593

$ws.Range("C25").Value = 95
$ws.Range("D25").Value = 64
$ws.Range("E25").Value = 48.4375
$ws.Range("F25").Value = 370
$ws.Range("G25").Value = 263
$ws.Range("H25").Value = 40.684410646387
$ws.Range("I25").Value = 1978
$ws.Range("J25").Value = 1611
$ws.Range("K25").Value = 22.780881440099
$ws.Range("L25").Value = 9.523809523809

$ws.Range("C26").Value = 126
$ws.Range("D26").Value = 106
$ws.Range("E26").Value = 18.867924528301
$ws.Range("F26").Value = 506
$ws.Range("G26").Value = 385
$ws.Range("H26").Value = 31.428571428571
$ws.Range("I26").Value = 2575
$ws.Range("J26").Value = 2297
$ws.Range("K26").Value = 12.102742707879
$ws.Range("L26").Value = 26.225490196078
$ws.Range("M26").Value = 8.329827513672

$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 25
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 92.307692307692
$ws.Range("I27").Value = 145
$ws.Range("J27").Value = 135
$ws.Range("K27").Value = 7.407407407407
$ws.Range("L27").Value = -14.201183431952

$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = -83.333333333333
$ws.Range("F28").Value = 37
$ws.Range("G28").Value = 26
$ws.Range("H28").Value = 42.307692307692
$ws.Range("I28").Value = 219
$ws.Range("J28").Value = 219
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -0.454545454545

$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 400
$ws.Range("F29").Value = 13
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 62.5
$ws.Range("I29").Value = 66
$ws.Range("J29").Value = 53
$ws.Range("K29").Value = 24.528301886792
$ws.Range("L29").Value = -31.25
$ws.Range("M29").Value = -37.735849056603
$ws.Range("N29").Value = -73.493975903614

$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 400
$ws.Range("F30").Value = 12
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 54
$ws.Range("J30").Value = 39
$ws.Range("K30").Value = 38.461538461538
$ws.Range("L30").Value = -26.027397260274
$ws.Range("M30").Value = -35.714285714285
$ws.Range("N30").Value = -76.521739130434

$ws.Range("D31").Value = 3
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = -75
$ws.Range("J31").Value = 37
$ws.Range("K31").Value = -56.756756756756

$ws.Range("D33").Value = 1
$ws.Range("G33").Value = 4
$ws.Range("J33").Value = 21
$ws.Range("K33").Value = -42.857142857142
$ws.Range("L33").Value = 0

